$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.853.08"
$ws.Range("E2").Value = "  -2.77%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.792.33"
$ws.Range("E3").Value = "  -0.55%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.27"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5320"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3843"
$ws.Range("E8").Value = "  +1.85%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07422"
$ws.Range("E9").Value = "  -0.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.33"
$ws.Range("E10").Value = "  -2.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.085"
$ws.Range("E11").Value = "  -2.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.177"
$ws.Range("E13").Value = "  +0.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.442"
$ws.Range("E14").Value = "  +1.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.28"
$ws.Range("E15").Value = "  -1.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.791.08"
$ws.Range("E16").Value = "  +1.33%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001059"
$ws.Range("E17").Value = "  -0.38%  "
$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "88.04"
$ws.Range("E18").Value = "  -2.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06515"
$ws.Range("E19").Value = "  +0.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9999"
$ws.Range("E20").Value = "  -0.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.23"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.946"
$ws.Range("E22").Value = "  +0.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.895.82"
$ws.Range("E23").Value = "  -2.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.13"
$ws.Range("E24").Value = "  +0.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.093"
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.00"
$ws.Range("E26").Value = "  -1.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.11"
$ws.Range("E27").Value = "  -1.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.997.46"
$ws.Range("E28").Value = "  -0.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.317"
$ws.Range("E29").Value = "  -1.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "121.18"
$ws.Range("E30").Value = "  -1.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1091"
$ws.Range("E31").Value = "  +4.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.099"
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.651"
$ws.Range("E33").Value = "  -1.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.496"
$ws.Range("E34").Value = "  -2.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06914"
$ws.Range("E35").Value = "  +7.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2201"
$ws.Range("E36").Value = "  -2.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02268"
$ws.Range("E37").Value = "  -1.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.032"
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.31"
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.362"
$ws.Range("E40").Value = "  -4.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6094"
$ws.Range("E41").Value = "  -1.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.171"
$ws.Range("E42").Value = "  -3.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.410"
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.29"
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.679"
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5688"
$ws.Range("E46").Value = "  -2.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "124.45"
$ws.Range("E47").Value = "  -0.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.909"
$ws.Range("E48").Value = "  -1.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.169"
$ws.Range("E49").Value = "  +1.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06790"
$ws.Range("E50").Value = "  -1.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.19"
$ws.Range("E51").Value = "  -1.01%  "
